# Milwaukee roster update: refresh player data (reorder, updated stats,
# drop George Hill / Jordan Nwora / Serge Ibaka, add Jae Crowder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the two trailing roster rows (18, 17) - bottom up so row
#    numbers of earlier rows are not disturbed while deleting.
# ---------------------------------------------------------------------
$ws.Range("A18:K18").EntireRow.Delete()
$ws.Range("A17:K17").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2. Drop all existing hyperlinks on the sheet. We will re-create the
#    15 that remain (K2:K16) further down, preserving their original
#    relationship order/targets.
# ---------------------------------------------------------------------
$ws.Range("A1:K16").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 3. Rewrite the player rows (2-16) with the refreshed roster data.
# ---------------------------------------------------------------------
$data = @(
  @{Row=2;  A=0;  B=5;  C="Jevon Carter";                  D="PG"; E="6-1";  F=200; G="September 14, 1995"; H="us"; I="4";  J="West Virginia";                 K="https://www.basketball-reference.com/players/c/carteje01.html"},
  @{Row=3;  A=1;  B=11; C="Brook Lopez";                    D="C";  E="7-0";  F=282; G="April 1, 1988";      H="us"; I="14"; J="Stanford";                      K="https://www.basketball-reference.com/players/l/lopezbr01.html"},
  @{Row=4;  A=2;  B=12; C="Grayson Allen";                  D="SG"; E="6-4";  F=198; G="October 8, 1995";    H="us"; I="4";  J="Duke";                          K="https://www.basketball-reference.com/players/a/allengr01.html"},
  @{Row=5;  A=3;  B=9;  C="Bobby Portis";                   D="PF"; E="6-10"; F=250; G="February 10, 1995";  H="us"; I="7";  J="Arkansas";                      K="https://www.basketball-reference.com/players/p/portibo01.html"},
  @{Row=6;  A=4;  B=34; C="Giannis Antetokounmpo";          D="PF"; E="7-0";  F=242; G="December 6, 1994";   H="gr"; I="9";  J=$null;                           K="https://www.basketball-reference.com/players/a/antetgi01.html"},
  @{Row=7;  A=5;  B=15; C="Jrue Holiday";                   D="PG"; E="6-3";  F=205; G="June 12, 1990";      H="us"; I="13"; J="UCLA";                          K="https://www.basketball-reference.com/players/h/holidjr01.html"},
  @{Row=8;  A=6;  B=24; C="Pat Connaughton";                D="SG"; E="6-5";  F=209; G="January 6, 1993";    H="us"; I="7";  J="Notre Dame";                    K="https://www.basketball-reference.com/players/c/connapa01.html"},
  @{Row=9;  A=7;  B=0;  C="MarJon Beauchamp";               D="SF"; E="6-6";  F=199; G="October 12, 2000";   H="us"; I="R";  J="Yakima Valley Community College"; K="https://www.basketball-reference.com/players/b/beaucma01.html"},
  @{Row=10; A=8;  B=23; C="Wesley Matthews";                D="SG"; E="6-4";  F=220; G="October 14, 1986";   H="us"; I="13"; J="Marquette";                     K="https://www.basketball-reference.com/players/m/matthwe02.html"},
  @{Row=11; A=9;  B=20; C="A.J. Green (TW)";                D="SG"; E="6-4";  F=200; G="September 27, 1999"; H="us"; I="R";  J="University of Northern Iowa";   K="https://www.basketball-reference.com/players/g/greenaj01.html"},
  @{Row=12; A=10; B=43; C="Thanasis Antetokounmpo";         D="SF"; E="6-6";  F=219; G="July 18, 1992";      H="gr"; I="4";  J=$null;                           K="https://www.basketball-reference.com/players/a/antetth01.html"},
  @{Row=13; A=11; B=7;  C="Joe Ingles";                     D="SF"; E="6-8";  F=220; G="October 2, 1987";    H="au"; I="8";  J=$null;                           K="https://www.basketball-reference.com/players/i/inglejo01.html"},
  @{Row=14; A=12; B=54; C="Sandro Mamukelashvili (TW)";     D="C";  E="6-11"; F=240; G="May 23, 1999";       H="us"; I="1";  J="Seton Hall";                    K="https://www.basketball-reference.com/players/m/mamuksa01.html"},
  @{Row=15; A=13; B=22; C="Khris Middleton";                D="SF"; E="6-7";  F=222; G="August 12, 1991";    H="us"; I="10"; J="Texas A&M";                     K="https://www.basketball-reference.com/players/m/middlkh01.html"},
  @{Row=16; A=14; B=$null; C="Jae Crowder";                 D="SF"; E="6-6";  F=235; G="July 6, 1990";       H="us"; I="10"; J="Marquette";                     K="https://www.basketball-reference.com/players/c/crowdja01.html"}
)

# Column I holds values that look numeric ("4", "14", ...) but must stay
# text (it also holds "R" for rookies). Force text formatting for that
# column up front so the written values keep their t="s" shared-string
# type instead of being auto-coerced to numbers.
$ws.Range("I2:I16").NumberFormat = "@"

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    if ($null -eq $rec.B) {
        $ws.Cells.Item($r, 2).ClearContents()
    } else {
        $ws.Cells.Item($r, 2).Value = $rec.B
    }
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
    $ws.Cells.Item($r, 8).Value = $rec.H
    $ws.Cells.Item($r, 9).Value = $rec.I
    if ($null -eq $rec.J) {
        $ws.Cells.Item($r, 10).ClearContents()
    } else {
        $ws.Cells.Item($r, 10).Value = $rec.J
    }
    $ws.Cells.Item($r, 11).Value = $rec.K
}

# Put column I formatting back to General now that the text values are
# committed, then restore the default (unstyled) cell style so no stray
# style index is left referenced on these cells.
$ws.Range("I2:I16").NumberFormat = "General"
$ws.Range("I2:I16").Style = "Normal"

# ---------------------------------------------------------------------
# 4. Re-create the hyperlinks for K2:K16, in order, so the relationship
#    ids line up the same way they did before the edit (rId1..rId15).
# ---------------------------------------------------------------------
foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 11), $rec.K) | Out-Null
    $ws.Cells.Item($r, 11).Style = "Hyperlink"
}
